$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "VERIFICATION EXPERIMENTS:" paragraph -> <w:ind w:hanging="0"/> becomes
#    <w:ind w:left="0" w:hanging="0"/>
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "VERIFICATION EXPERIMENTS:*") {
        $p.Format.LeftIndent = 0
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "**" + "Modificare thoughput" -> merge into a single run "**Modificare
#    thoughput" and set the same indent fix on that paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Modificare thoughput*") {
        $p.Format.LeftIndent = 0
        break
    }
}
$d.Content.Find.Execute("**Modificare thoughput", $true, $false, $false, $false, $false, $true, 1, $false, "**Modificare thoughput", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "**Aggiungere controllo ..." paragraph -> same indent fix, no text change
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Aggiungere controllo*") {
        $p.Format.LeftIndent = 0
        break
    }
}

# ---------------------------------------------------------------------------
# 4) "Un coppia tx-rx, un canale, ..." -> insert "(> timeslot)" before the
#    ", buffer 5 slot" tail, as its own run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Un coppia tx-rx, un canale, Probabilità di trasmissione a 1, Interarrival time deterministico, buffer 5 slot", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Un coppia tx-rx, un canale, Probabilità di trasmissione a 1, Interarrival time deterministico (> timeslot), buffer 5 slot", `
    2) | Out-Null

$rng4 = $d.Content
$rng4.Find.Execute("(> timeslot)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s4 = $rng4.Start
$e4 = $rng4.End
$mid4 = $d.Range($s4, $e4)
$mid4.Bold = 0
$mid4.Bold = 1

# ---------------------------------------------------------------------------
# 5) "Una coppia tx-rx, 5 canali, ..." -> insert ": < timeslot" right before
#    the closing ")" as its own run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Una coppia tx-rx, 5 canali, Probabilità di trasmissione a 1, interarrival time deterministico (con pacchetto sempre nel buffer), buffer 5 slot", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Una coppia tx-rx, 5 canali, Probabilità di trasmissione a 1, interarrival time deterministico (con pacchetto sempre nel buffer : < timeslot), buffer 5 slot", `
    2) | Out-Null

$rng5 = $d.Content
$rng5.Find.Execute(": < timeslot", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s5 = $rng5.Start
$e5 = $rng5.End
$mid5 = $d.Range($s5, $e5)
$mid5.Bold = 0
$mid5.Bold = 1

# ---------------------------------------------------------------------------
# 6) "Un coppia, un canale, probabilità ..." -> insert "(> timeslot)" before
#    the ", buffer 5 slot" tail, as its own run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Un coppia, un canale, probabilità di trasmissione a 0.5, interarrival time deterministico, buffer 5 slot", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Un coppia, un canale, probabilità di trasmissione a 0.5, interarrival time deterministico (> timeslot), buffer 5 slot", `
    2) | Out-Null

$rng6 = $d.Range($e5, $d.Content.End)
$rng6.Find.Execute("(> timeslot)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s6 = $rng6.Start
$e6 = $rng6.End
$mid6 = $d.Range($s6, $e6)
$mid6.Bold = 0
$mid6.Bold = 1

# ---------------------------------------------------------------------------
# 7) "Degeneracy test (Probabilità trasmissione a " + "0, 0 coppie tx-rx" ->
#    merge into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Degeneracy test (Probabilità trasmissione a 0, 0 coppie tx-rx", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Degeneracy test (Probabilità trasmissione a 0, 0 coppie tx-rx", `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 8) "C" + "onsistency test (...)" -> merge into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Consistency test (1 coppia, 1 canale, interarrival deterministico a k vs 2 coppie 500 canali, interarrival deterministico a 2k)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Consistency test (1 coppia, 1 canale, interarrival deterministico a k vs 2 coppie 500 canali, interarrival deterministico a 2k)", `
    2) | Out-Null
